$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted as row 285 ("Feria Lagunitas
# de Puerto Montt" - Choclo - Choclero - Primera, dated 2023-04-11),
# pushing every existing row from 285 downward by one (old row 285 -> new
# row 286, ..., old row 382 -> new row 383), growing the used range from
# A1:R382 to A1:R383.
$ws.Rows.Item(285).Insert()

$ws.Cells.Item(285, 1).Value  = 4
$ws.Cells.Item(285, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(285, 3).Value  = "Los Lagos"
$ws.Cells.Item(285, 4).Value  = 45027
$ws.Cells.Item(285, 5).Value  = 10
$ws.Cells.Item(285, 6).Value  = 100112024
$ws.Cells.Item(285, 7).Value  = "Choclo"
$ws.Cells.Item(285, 8).Value  = "Choclero"
$ws.Cells.Item(285, 9).Value  = "Primera"
$ws.Cells.Item(285, 10).Value = 6000
$ws.Cells.Item(285, 11).Value = 550
$ws.Cells.Item(285, 12).Value = 600
$ws.Cells.Item(285, 13).Value = 575
$ws.Cells.Item(285, 14).Value = "$/unidad"
$ws.Cells.Item(285, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(285, 16).Value = 575
$ws.Cells.Item(285, 17).Value = 1
$ws.Cells.Item(285, 18).Value = "Hortaliza"
